# channel config spreadsheet update
#
# - swap the SV-HE-201/202 "state" and "position" channel IDs
#   (B7<->B8, B9<->B10 on the "channels" sheet)
# - make "channels" the active sheet/tab again (was "calibration"),
#   with the cursor left on B11
# - "calibration" is no longer the selected tab

$wb = $excel.ActiveWorkbook

$channels = $wb.Worksheets.Item("channels")

# --- swap the ID values in B7:B10 -----------------------------------------
$b7 = $channels.Range("B7").Value2
$b8 = $channels.Range("B8").Value2
$channels.Range("B7").Value = $b8
$channels.Range("B8").Value = $b7

$b9 = $channels.Range("B9").Value2
$b10 = $channels.Range("B10").Value2
$channels.Range("B9").Value = $b10
$channels.Range("B10").Value = $b9

# --- restore "channels" as the active sheet/tab with cursor at B11 --------
$channels.Activate() | Out-Null
$channels.Range("B11").Select() | Out-Null
